# ===========================================================================
# Memoli.xlsx -- "Memoli 2015 + 2016 model fits"
#
#   * Renames the "2015" sheet to "Sheet1"
#   * Relabels the existing 2015 "Mean" rows (A2:A11) as "2015Mean"
#   * Adds five serum-cytokine columns (GCSF, IFNG, IL6, TNFA, TSS) in D:H,
#     with sparse values only for a handful of rows
#   * Appends two new volunteer-group blocks for the 2016 HAI study
#     (2016LoHAI -> rows 12-21, 2016HiHAI -> rows 22-28)
#   * Re-creates the original alternating number-format / fill styling on
#     the new cells by copying formats from existing same-styled cells
#     (so the style table itself is reused rather than duplicated)
# ===========================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet "2015" -> "Sheet1" -----------------------------------
$ws.Name = "Sheet1"

# --- Relabel existing 2015 Mean rows (A2:A11): "Mean" -> "2015Mean" ----
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("A$r").Value = "2015Mean"
}

# --- New cytokine header row (D1:H1) -----------------------------------
$headers = @{ "D1" = "GCSF"; "E1" = "IFNG"; "F1" = "IL6"; "G1" = "TNFA"; "H1" = "TSS" }
foreach ($ref in $headers.Keys) {
    $ws.Range($ref).Value = $headers[$ref]
}

# --- Sparse cytokine values for the 2015 Mean rows (D2:H11) ------------
# Row -> @{ col = value }
$cytokines = @{
    2  = @{ D = 2.07214710083425;    E = 1.8787231309594017;  F = 0.83520630659662332; G = 1.3160853443758567;  H = 0.114293580562431 }
    3  = @{ D = 2.0925756194889131;  E = 1.9180116497471442;  F = 0.97285842091697461; G = 1.3224607408897178;  H = 0.58277995074171596 }
    4  = @{                                                                                                      H = 1.8232716956400701 }
    5  = @{ D = 2.1704207543291636;  E = 1.7665154267369445;  F = 0.99234476666539417; G = 1.5689442921407151;  H = 3.49150548993437 }
    6  = @{                                                                                                      H = 3.6574364744715799 }
    7  = @{ D = 1.8875998384572792;  E = 1.6309935214254174;  F = 0.99124375149873178; G = 1.4213580892700912;  H = 3.42700586849893 }
    8  = @{                                                                                                      H = 2.1115822801856599 }
    9  = @{                                                                                                      H = 1.8289696160786599 }
    10 = @{                                                                                                      H = 0.67005030569462798 }
    11 = @{                                                                                                      H = 0.60651430315689803 }
}
foreach ($row in $cytokines.Keys) {
    $cols = $cytokines[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# --- New volunteer-group blocks: A/B/C for rows 12-28 -------------------
# Row -> (DAY, VALUE); rows 12-21 belong to 2016LoHAI, rows 22-28 to 2016HiHAI
$lohai = @(
    1122.948171278714, 4515.192378428178, 14397.045636220895, 19398.573030675019,
    9053.8757952982342, 1240.2955338330044, 541.7675012235203, 5417.6750122352041,
    43.680516977211418, 0
)
$hihai = @(
    0, 248.704391122673, 6083.7682475720385, 139.27876015629599,
    187.6641408787614, 13.927876015629595, 0
)

$row = 12
for ($day = 1; $day -le $lohai.Count; $day++) {
    $ws.Range("A$row").Value = "2016LoHAI"
    $ws.Range("B$row").Value = $day
    $ws.Range("C$row").Value = $lohai[$day - 1]
    $row++
}
for ($day = 1; $day -le $hihai.Count; $day++) {
    $ws.Range("A$row").Value = "2016HiHAI"
    $ws.Range("B$row").Value = $day
    $ws.Range("C$row").Value = $hihai[$day - 1]
    $row++
}

# --- Re-apply the workbook's existing alternating-row styles to the new
#     cells, by copying formats from a donor cell that already carries
#     that exact style (keeps the shared style table from growing) ------
$xlPasteFormats = -4122

# Donor cells per style, taken from the original, untouched B/C columns:
#   style "1" (0.00 + light fill) -> C2
#   style "2" (light fill only)   -> B2
#   style "3" (0.00 only)         -> C3
$style1 = "C2"
$style2 = "B2"
$style3 = "C3"

function Copy-FormatOnly($donorRef, $targetRef) {
    $ws.Range($donorRef).Copy()
    $ws.Range($targetRef).PasteSpecial($xlPasteFormats)
}

# D2:G2 / D3:G3 / D5:G5 / D7:G7 alternate style3/style1/style1/style3,
# mirroring the C-column pattern already present for those same rows.
$cytokineStyle = @{
    2 = @{ D = $style3; E = $style1; F = $style1; G = $style3 }
    3 = @{ D = $style1; E = $style3; F = $style3; G = $style1 }
    5 = @{ D = $style3; E = $style1; F = $style1; G = $style3 }
    7 = @{ D = $style3; E = $style1; F = $style1; G = $style3 }
}
foreach ($row in $cytokineStyle.Keys) {
    $cols = $cytokineStyle[$row]
    foreach ($col in $cols.Keys) {
        Copy-FormatOnly $cols[$col] "$col$row"
    }
}

# 2016LoHAI / 2016HiHAI rows that carry the light-fill "style 2" on B & C
# (every other row, matching the source workbook's banding)
$style2Rows = 13, 15, 16, 18, 19, 21, 23, 24, 28
foreach ($row in $style2Rows) {
    Copy-FormatOnly $style2 "B$row"
    Copy-FormatOnly $style2 "C$row"
}

$excel.CutCopyMode = $false

# --- Match the final selection shown in the diff (activeCell F17) ------
$ws.Range("F17").Select() | Out-Null
